# Updates to Batching and Profiles
# Accuracy and Evasion are now agent properties with ranges.
#
# Insert two new columns (B, C) for the new agent properties "Accuracy" and
# "Evasion", ahead of the existing property columns (which shift right).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns before column B; existing Curiosity..Efficiency
# columns (B:H) shift right to D:J.
$ws.Columns("B:C").Insert()

# Headers for the new columns.
$ws.Range("B1").Value = "Accuracy"
$ws.Range("C1").Value = "Evasion"

# Default data rows: numeric 0 baseline (rows 2-3) and the "…" placeholder
# row (row 4), matching the pattern already used by the other property
# columns.
$ws.Range("B2:C3").Value = 0
$ws.Range("B4:C4").Value = "…"

# Leave the selection on the new first inserted data cell, matching the
# author's final cursor position.
$ws.Range("D4").Select() | Out-Null
